# DEV 12 - Project Global Deletion
#
# For the selected project, every application that was in the
# "Withdrawal Pending" state is moved to "Deleted" (global deletion of a
# project cascades into deleting its pending-withdrawal applications),
# and the Date column is refreshed to the time of the action.
# A brand-new application row is then appended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark previously "Withdrawal Pending" applications as "Deleted" ---
$ws.Range("D2").Value = "Deleted"
$ws.Range("F2").Value = 45771.12818586805

$ws.Range("D3").Value = "Deleted"
$ws.Range("F3").Value = 45771.128186712966

# --- Append the new application record in row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "T7654321B"
$ws.Range("D5").Value = "Pending"
$ws.Range("E5").Value = "2-ROOM"
$ws.Range("F5").Value = 45769.900285069445

# Columns A/B carry a legacy 2-decimal column style; match the rest of the
# table (rows 2-4) by keeping the new row on the default "Normal" style.
$ws.Range("A5:B5").Style = "Normal"

# --- Update the sheet view to reflect the newly added row ---
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("A5:F5").Select()
